$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new task row "Оптимизация ИИ" below the last existing row (row 17)
$ws.Range("A18").Value = "Оптимизация ИИ"

# Apply red fill style (same as used on B7/B8, style index 1) to B17 and B18
$ws.Range("B17").Interior.Color = $ws.Range("B7").Interior.Color
$ws.Range("B18").Interior.Color = $ws.Range("B7").Interior.Color

# Update selection to match the target state
$ws.Range("D16").Select()
